$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F (reuse the same header style as the other header cells)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 13:41:09.103204",
    "2021-10-05 13:41:09.103214",
    "2021-10-05 13:41:09.103216",
    "2021-10-05 13:41:09.103218",
    "2021-10-05 13:41:09.103221",
    "2021-10-05 13:41:09.103223",
    "2021-10-05 13:41:09.103224",
    "2021-10-05 13:41:09.103226",
    "2021-10-05 13:41:09.103228",
    "2021-10-05 13:41:09.103230",
    "2021-10-05 13:41:09.103232",
    "2021-10-05 13:41:09.103234",
    "2021-10-05 13:41:09.103236",
    "2021-10-05 13:41:09.103238",
    "2021-10-05 13:41:09.103240",
    "2021-10-05 13:41:09.103242",
    "2021-10-05 13:41:09.103244",
    "2021-10-05 13:41:09.103246",
    "2021-10-05 13:41:09.103248",
    "2021-10-05 13:41:09.103250",
    "2021-10-05 13:41:09.103252",
    "2021-10-05 13:41:09.103254",
    "2021-10-05 13:41:09.103255",
    "2021-10-05 13:41:09.103257",
    "2021-10-05 13:41:09.103260",
    "2021-10-05 13:41:09.103261",
    "2021-10-05 13:41:09.103263",
    "2021-10-05 13:41:09.103265",
    "2021-10-05 13:41:09.103267",
    "2021-10-05 13:41:09.103269",
    "2021-10-05 13:41:09.103271",
    "2021-10-05 13:41:09.103273",
    "2021-10-05 13:41:09.103275",
    "2021-10-05 13:41:09.103278",
    "2021-10-05 13:41:09.103280",
    "2021-10-05 13:41:09.103282",
    "2021-10-05 13:41:09.103284",
    "2021-10-05 13:41:09.103286",
    "2021-10-05 13:41:09.103288",
    "2021-10-05 13:41:09.103290",
    "2021-10-05 13:41:09.103292",
    "2021-10-05 13:41:09.103294",
    "2021-10-05 13:41:09.103296",
    "2021-10-05 13:41:09.103298",
    "2021-10-05 13:41:09.103300",
    "2021-10-05 13:41:09.103302",
    "2021-10-05 13:41:09.103304",
    "2021-10-05 13:41:09.103307"
)

$row = 2
foreach ($ts in $timestamps) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $ts
    $row = $row + 1
}
